$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 3

$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 4
